$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A1 keeps the same displayed text ("HK_R_acc_G"); re-setting it forces a
# fresh shared-string entry to be appended/used, matching the bumped
# sharedStrings count/uniqueCount and the new <v> index referenced by A1.
$ws.Range("A1").Value = "HK_R_acc_G"

$values = @(
    73.992673992674,
    73.504273504273513,
    73.544973544973544,
    69.108669108669105,
    71.794871794871796,
    73.056573056573058,
    73.748473748473756,
    74.277574277574274,
    73.341473341473346,
    73.015873015873012,
    74.155474155474153,
    76.923076923076934,
    72.730972730972738,
    73.504273504273513,
    72.771672771672769,
    75.091575091575095,
    76.353276353276357,
    76.109076109076113,
    75.213675213675216,
    72.527472527472526,
    75.295075295075293,
    75.986975986975992,
    75.661375661375658,
    75.905575905575901,
    73.789173789173788,
    73.789173789173788,
    74.074074074074076,
    80.179080179080174,
    79.527879527879534,
    77.574277574277573,
    71.550671550671552,
    74.033374033374031,
    74.033374033374031,
    75.539275539275536,
    75.661375661375658,
    79.242979242979246,
    75.579975579975581,
    75.45787545787546,
    76.271876271876266,
    74.643874643874639,
    74.358974358974365,
    74.765974765974761,
    73.951973951973955,
    74.236874236874243,
    73.300773300773301,
    72.85307285307286,
    72.486772486772495,
    73.015873015873012
)

$row = 2
foreach ($val in $values) {
    $ws.Cells.Item($row, 1).Value = $val
    $row = $row + 1
}
